$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 3360
$ws.Range("I86").Value = 2500
$ws.Range("J86").Value = 3575
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 3575
$ws.Range("M86").Value = -1377
$ws.Range("N86").Value = -5821

# Row 89
$ws.Range("H89").Value = 3360
$ws.Range("I89").Value = 2500
$ws.Range("J89").Value = 3575
$ws.Range("K89").Value = 12500
$ws.Range("L89").Value = 17875
$ws.Range("M89").Value = -6884
$ws.Range("N89").Value = -29107

# Row 132
$ws.Range("H132").Value = 4447015
$ws.Range("I132").Value = 5265462.5
$ws.Range("K132").Value = 15796387.5
$ws.Range("M132").Value = -15793857.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4509.98
$ws.Range("I32").Value = 4187.7656
$ws.Range("J32").Value = 5883.6313
$ws.Range("K32").Value = 4187.7656
$ws.Range("L32").Value = 5883.6313
$ws.Range("M32").Value = -3900.7656
$ws.Range("N32").Value = -6457.6313

# Row 41
$ws.Range("H41").Value = 18636.5
$ws.Range("I41").Value = 1494.6666
$ws.Range("J41").Value = 70062
$ws.Range("K41").Value = 1494.6666
$ws.Range("L41").Value = 70062
$ws.Range("M41").Value = -1080.6666
$ws.Range("N41").Value = -70890

# Row 63
$ws.Range("H63").Value = 5516.5
$ws.Range("I63").Value = 3100
$ws.Range("J63").Value = 6724.75
$ws.Range("K63").Value = 3100
$ws.Range("L63").Value = 6724.75
$ws.Range("M63").Value = -2414
$ws.Range("N63").Value = -8096.75

# Row 66
$ws.Range("H66").Value = 5516.5
$ws.Range("I66").Value = 3100
$ws.Range("J66").Value = 6724.75
$ws.Range("K66").Value = 15500
$ws.Range("L66").Value = 33623.75
$ws.Range("M66").Value = -12068
$ws.Range("N66").Value = -40487.75

# Row 132
$ws.Range("H132").Value = 1846.4603
$ws.Range("I132").Value = 1298.2444
$ws.Range("J132").Value = 3217
$ws.Range("K132").Value = 3894.7332
$ws.Range("L132").Value = 9651
$ws.Range("M132").Value = -1364.7332
$ws.Range("N132").Value = -14711

# Row 139
$ws.Range("H139").Value = 26800
$ws.Range("J139").Value = 26800
$ws.Range("L139").Value = 26800
$ws.Range("N139").Value = -37080

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 26411.9
$ws.Range("J2").Value = 32626.875
$ws.Range("L2").Value = 32626.875
$ws.Range("N2").Value = -32852.875

# Row 37
$ws.Range("H37").Value = 20000
$ws.Range("J37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("N37").Value = -20214

# Row 45
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

# Row 132
$ws.Range("H132").Value = 1598.2778
$ws.Range("I132").Value = 1178.6154
$ws.Range("J132").Value = 2689.4
$ws.Range("K132").Value = 3535.8462
$ws.Range("L132").Value = 8068.200000000001
$ws.Range("M132").Value = -1005.8462
$ws.Range("N132").Value = -13128.2

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1191.6471
$ws.Range("I5").Value = 604.1429000000001
$ws.Range("K5").Value = 1812.4287
$ws.Range("M5").Value = -1700.4287

# Row 131
$ws.Range("H131").Value = 1314.421
$ws.Range("I131").Value = 1699.4117
$ws.Range("J131").Value = 1150.8
$ws.Range("K131").Value = 5098.2351
$ws.Range("L131").Value = 3452.4
$ws.Range("M131").Value = -58.23509999999987
$ws.Range("N131").Value = -13532.4

# Row 135
$ws.Range("H135").Value = 1191.6471
$ws.Range("I135").Value = 604.1429000000001
$ws.Range("K135").Value = 5437.2861
$ws.Range("M135").Value = -2902.2861

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3632.2104
$ws.Range("I126").Value = 2502.4
$ws.Range("J126").Value = 4035.7144
$ws.Range("K126").Value = 7507.200000000001
$ws.Range("L126").Value = 12107.1432
$ws.Range("M126").Value = -5037.200000000001
$ws.Range("N126").Value = -17047.1432

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 224848.28
$ws.Range("J2").Value = 41413.8
$ws.Range("L2").Value = 41413.8
$ws.Range("N2").Value = -41637.8

# Row 35
$ws.Range("H35").Value = 21000
$ws.Range("I35").Value = 3000
$ws.Range("K35").Value = 3000
$ws.Range("M35").Value = -2664

# Row 55
$ws.Range("H55").Value = 1467.2
$ws.Range("I55").Value = 500
$ws.Range("K55").Value = 500
$ws.Range("M55").Value = -327

# Row 132
$ws.Range("H132").Value = 1739.509
$ws.Range("I132").Value = 1020.3488
$ws.Range("K132").Value = 3061.0464
$ws.Range("M132").Value = -531.0464000000002

# Row 136
$ws.Range("H136").Value = 1785.6809
$ws.Range("I136").Value = 1357.1538
$ws.Range("K136").Value = 4071.4614
$ws.Range("M136").Value = -1521.4614

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 5845001
$ws.Range("J2").Value = 5014001
$ws.Range("L2").Value = 5014001
$ws.Range("N2").Value = -5014225

# Row 3
$ws.Range("H3").Value = 2510145
$ws.Range("I3").Value = 3333526.8
$ws.Range("J3").Value = 40000
$ws.Range("K3").Value = 3333526.8
$ws.Range("L3").Value = 40000
$ws.Range("M3").Value = -3333412.8
$ws.Range("N3").Value = -40228

# Row 4
$ws.Range("H4").Value = 1795939.1
$ws.Range("I4").Value = 22000
$ws.Range("J4").Value = 2781461
$ws.Range("K4").Value = 22000
$ws.Range("L4").Value = 2781461
$ws.Range("M4").Value = -21887
$ws.Range("N4").Value = -2781687

# Row 6
$ws.Range("H6").Value = 14564.2
$ws.Range("I6").Value = 504.5
$ws.Range("J6").Value = 23937.334
$ws.Range("K6").Value = 504.5
$ws.Range("L6").Value = 23937.334
$ws.Range("M6").Value = -389.5
$ws.Range("N6").Value = -24167.334

# Row 132
$ws.Range("H132").Value = 12184.412
$ws.Range("I132").Value = 2195.575
$ws.Range("J132").Value = 48507.453
$ws.Range("K132").Value = 6586.724999999999
$ws.Range("L132").Value = 145522.359
$ws.Range("M132").Value = -4056.724999999999
$ws.Range("N132").Value = -150582.359

# Row 135
$ws.Range("H135").Value = 39005
$ws.Range("J135").Value = 39005
$ws.Range("L135").Value = 39005
$ws.Range("N135").Value = -49145

# Row 136
$ws.Range("H136").Value = 992.98
$ws.Range("I136").Value = 689.7368
$ws.Range("K136").Value = 2069.2104
$ws.Range("M136").Value = 480.7896000000001

# Row 141
$ws.Range("H141").Value = 28423.076
$ws.Range("J141").Value = 28423.076
$ws.Range("L141").Value = 28423.076
$ws.Range("N141").Value = -38783.076
